# Auto-generated edit script: applies the numeric "want-to-go" counter
# bumps and the rolling event-list refresh described in the commit diff.
$wb = $excel.ActiveWorkbook

# --- "展览" (Exhibition) sheet: F-column counter bumps ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 225
$ws1.Range("F3").Value = 54947
$ws1.Range("F4").Value = 3114
$ws1.Range("F5").Value = 5268
$ws1.Range("F6").Value = 1177
$ws1.Range("F10").Value = 1098
$ws1.Range("F11").Value = 1363
$ws1.Range("F12").Value = 119
$ws1.Range("F14").Value = 219
$ws1.Range("F15").Value = 396
$ws1.Range("F16").Value = 55
$ws1.Range("F21").Value = 5396
$ws1.Range("F23").Value = 5284
$ws1.Range("F24").Value = 9297
$ws1.Range("F27").Value = 149
$ws1.Range("F28").Value = 239
$ws1.Range("F29").Value = 451
$ws1.Range("F31").Value = 106
$ws1.Range("F32").Value = 4274
$ws1.Range("F33").Value = 286

# --- "演出" (Performance) sheet: F-column counter bump ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 139

# --- "本地生活" (Local Life) sheet: F-column counter bump ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 592

# --- "全部类型" (All types) sheet: F-column counter bumps ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 592
$ws4.Range("F4").Value = 225
$ws4.Range("F5").Value = 3114
$ws4.Range("F7").Value = 139
$ws4.Range("F8").Value = 1177
$ws4.Range("F13").Value = 1098
$ws4.Range("F15").Value = 1363
$ws4.Range("F17").Value = 119
$ws4.Range("F18").Value = 219
$ws4.Range("F20").Value = 396
$ws4.Range("F21").Value = 55
$ws4.Range("F26").Value = 5396
$ws4.Range("F28").Value = 5284
$ws4.Range("F29").Value = 9297
$ws4.Range("F33").Value = 149
$ws4.Range("F34").Value = 239
$ws4.Range("F35").Value = 451
$ws4.Range("F40").Value = 107
$ws4.Range("F41").Value = 4274
$ws4.Range("F47").Value = 286

# --- "全部类型": rows 43-49 event list refresh (new shows pushed in) ---
$ws4.Range("B43").Value = "'2024-12-22"
$ws4.Range("C43").Value = "杭州·《你的名字》《天气之子》《铃芽之旅》——新海诚动漫三部曲钢琴演奏会"
$ws4.Range("E43").Value = "2024.12.22 19:30-12.22 21:00"
$ws4.Range("F43").Value = 90
$ws4.Range("G43").Value = 153
$ws4.Range("H43").Value = "https://show.bilibili.com/platform/detail.html?id=88648"
$ws4.Range("I43").Value = "//i1.hdslb.com/bfs/openplatform/202407/nEB3TPxP1720064877363.jpeg"
$ws4.Range("B44").Value = "'2024-12-24"
$ws4.Range("C44").Value = "杭州·百老汇悬疑惊悚喜剧《死亡陷阱》"
$ws4.Range("D44").Value = "武林广场29号 杭州剧院"
$ws4.Range("E44").Value = "2024.12.24 19:30-12.25 22:35"
$ws4.Range("F44").Value = 14
$ws4.Range("G44").Value = 280
$ws4.Range("H44").Value = "https://show.bilibili.com/platform/detail.html?id=88084"
$ws4.Range("I44").Value = "//i1.hdslb.com/bfs/openplatform/202407/KIysvIdZ1719818362132.jpeg"
$ws4.Range("C45").Value = "杭州·维也纳皇家交响乐团2025新年音乐会"
$ws4.Range("D45").Value = "桥弄街399号 杭州运河大剧院"
$ws4.Range("E45").Value = "2024.12.24 19:30-12.24 21:15"
$ws4.Range("H45").Value = "https://show.bilibili.com/platform/detail.html?id=91492"
$ws4.Range("I45").Value = "//i2.hdslb.com/bfs/openplatform/202408/bGIVRHpJ1724661705028.jpeg"
$ws4.Range("B46").Value = "'2024-12-31"
$ws4.Range("C46").Value = "杭州·2025大剧院缤纷跨年夜 爱·大声告白-成都“知更”室内合唱团音乐会"
$ws4.Range("D46").Value = "新业路39号 杭州大剧院"
$ws4.Range("E46").Value = "2024.12.31 22:30-2025.01.01 00:00"
$ws4.Range("F46").Value = 3
$ws4.Range("H46").Value = "https://show.bilibili.com/platform/detail.html?id=91210"
$ws4.Range("I46").Value = "//i0.hdslb.com/bfs/openplatform/202408/RGm2uKFJ1724395472501.jpeg"
$ws4.Range("F47").Value = 286
$ws4.Range("B48").Value = "'2024-12-31"
$ws4.Range("C48").Value = "杭州·法国爱乐乐团2025新年音乐会"
$ws4.Range("D48").Value = "桥弄街399号 杭州运河大剧院"
$ws4.Range("E48").Value = "2024.12.31 19:30-12.31 21:15"
$ws4.Range("F48").Value = 9
$ws4.Range("G48").Value = 280
$ws4.Range("H48").Value = "https://show.bilibili.com/platform/detail.html?id=91493"
$ws4.Range("I48").Value = "//i2.hdslb.com/bfs/openplatform/202408/8IHtZOW41724664782598.jpeg"
$ws4.Range("B49").Value = "'2025-01-01"
$ws4.Range("C49").Value = "杭州·【早鸟优惠】大型正版授权互动卡通儿童剧《海底小纵队之深海探秘》"
$ws4.Range("D49").Value = "湖墅南路136-138号 浙话艺术剧院"
$ws4.Range("E49").Value = "2025.01.01 10:30-01.01 11:40"
$ws4.Range("F49").Value = 1
$ws4.Range("G49").Value = 40
$ws4.Range("H49").Value = "https://show.bilibili.com/platform/detail.html?id=92951"
$ws4.Range("I49").Value = "//i2.hdslb.com/bfs/openplatform/202409/oZlaKX931727335820196.jpeg"
